# Appends the "Slots" / "Formularios" notes section (26/03/2024) to the end of
# the document, right after the "Las acciones customizadas..." paragraph.
# Built as WordprocessingML (OOXML) paragraph fragments and inserted in one
# shot via Range.InsertXML so formatting (list numbering, page break, etc.)
# matches exactly.

$d = $word.ActiveDocument

# Each element below is one <w:p>...</w:p> paragraph to append, in order.
$newParagraphs = @(
    '<w:p><w:r><w:br w:type="page"/></w:r></w:p>',
    '<w:p><w:r><w:lastRenderedPageBreak/><w:t>26/03/2024</w:t></w:r></w:p>',
    '<w:p><w:r><w:t>Slots</w:t></w:r></w:p>',
    '<w:p><w:r><w:t>Concepto de slots o ranuras: mecanismo para almacenar ciertos datos durante la conversación. Esta información se gurda durante el transcurso de la información y se recuperan en un momento determinado mas adelante.</w:t></w:r></w:p>',
    '<w:p><w:r><w:t>Ejemplo:</w:t></w:r></w:p>',
    '<w:p><w:r><w:t>Si un usuario indica la talla de la ropa, esta talla se guarda para luego utilizarla mas adelante. Se puede definir en el archivo nlu.yml con varios intents donde se reconozcan entidades. En el archivo de stories.yml se declara</w:t></w:r><w:r><w:t xml:space="preserve"> la variable donde se va a almacenar. En el domain.yml hay que declarar también las entities. También debe existir una sección slots con la variable. Finalmenete en el código de la acción se puede utilizar el valor almacenado en el slot</w:t></w:r></w:p>',
    '<w:p/>',
    '<w:p><w:r><w:t>Formularios</w:t></w:r></w:p>',
    '<w:p><w:r><w:t>Permite manejar varios datos al mismo tiempo.</w:t></w:r></w:p>',
    '<w:p><w:r><w:t>Pasos a seguir</w:t></w:r></w:p>',
    '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:t>Agregar ‘forms’ al domain.yml</w:t></w:r></w:p>',
    '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:t>Paso 2</w:t></w:r></w:p>',
    '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:t>Paso 3: asociar regla en rules.yml</w:t></w:r></w:p>',
    '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:t>Paso 4: activar formulario en stories.yml</w:t></w:r></w:p>',
    '<w:p><w:r><w:t>Durante el transcurso de un formulario también pueden ocurrir situaciones que desvíen el flujo del formulario:</w:t></w:r></w:p>',
    '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:t>Referencias a otras intencioes</w:t></w:r></w:p>',
    '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:t>Peticion de cancelación</w:t></w:r></w:p>',
    '<w:p><w:r><w:t>Referencia a otras intenciones: utilizar reglas</w:t></w:r></w:p>',
    '<w:p><w:r><w:t>Peticion de cancelación del formulario: requiere de los siguientes pasos</w:t></w:r></w:p>',
    '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:t>Respuesta de cancelacion</w:t></w:r></w:p>'
)

$newBodyXml = [string]::Join("", $newParagraphs)

$xmlDoc = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
    '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData>' +
    '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
    '<w:body>' + $newBodyXml + '</w:body></w:document>' +
    '</pkg:xmlData></pkg:part></pkg:package>'

# Anchor the insertion at the very end of the document content (after the
# last paragraph, "Las acciones customizadas se definen en el archivo
# actions.py", and before the final section properties) using an explicit
# Range built from Content.End - collapsing a Range obtained from
# Paragraphs.Last can clobber the preceding paragraph's text in this runtime,
# so we avoid that pattern.
$endPos = $d.Content.End
$rng = $d.Range($endPos, $endPos)
$rng.InsertXML($xmlDoc)

Write-Host "Inserted" $newParagraphs.Count "paragraphs. Document now has" $d.Paragraphs.Count "paragraphs."
